$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transcriptions")

# ---------------------------------------------------------------------------
# 1. Insert 11 new rows right below M54 (old row 24) for its sub-manifests,
#    M54-1 through M54-11. This also naturally shifts every row below
#    (old 25-37) down by 11, carrying the SUM formula / notes along with it.
# ---------------------------------------------------------------------------
$ws.Range("A25:A35").EntireRow.Insert()

$subNumbers = 1..11
$startRow = 25
for ($i = 0; $i -lt $subNumbers.Count; $i++) {
    $r = $startRow + $i
    $row = $ws.Rows.Item($r)
    $row.RowHeight = 16

    $ws.Range("A" + $r).Style = "Normal"
    $ws.Range("A" + $r).WrapText = $true
    $ws.Range("A" + $r).Value = "M54-" + $subNumbers[$i]

    # B/C/D/F/G inherited formatting from the row above on insert - clear it
    # back down to the (unstyled) blank cells the new rows actually have.
    $ws.Range("B" + $r + ":D" + $r).ClearContents()
    $ws.Range("B" + $r + ":D" + $r).Style = "Normal"
    $ws.Range("F" + $r + ":G" + $r).ClearContents()
    $ws.Range("F" + $r + ":G" + $r).Style = "Normal"

    # E column keeps the (empty) Hyperlink-style placeholder seen elsewhere
    # in the sheet for rows awaiting a manifest link.
    $ws.Range("E" + $r).ClearContents()
    $ws.Range("E" + $r).Style = "Hyperlink"
    $ws.Range("E" + $r).WrapText = $true
}

# ---------------------------------------------------------------------------
# 2. Row 24 (M54): the page-count (F24) no longer applies to the parent row
#    now that it has been split into sub-manifests, and two new "N" flags
#    (Transcription-Facsimile Links? / Transcription Redraft?) are recorded
#    in H24/I24 using Excel's built-in "Bad" (red) cell style.
# ---------------------------------------------------------------------------
$ws.Range("F24").ClearContents()
$ws.Range("F24").Style = "Normal"

$ws.Range("H24:I24").Style = "Bad"
$ws.Range("H24:I24").Font.Size = 12
$ws.Range("H24:I24").HorizontalAlignment = -4108
$ws.Range("H24:I24").VerticalAlignment = -4108
$ws.Range("H24:I24").WrapText = $true
$ws.Range("H24:I24").Value = "N"

# ---------------------------------------------------------------------------
# 3. The three rows that used to be 25-27 (now 36-38 = M108/M109/M110)
#    gain a "Y" in the new Transcription-Facsimile Links? column (H).
# ---------------------------------------------------------------------------
$ws.Range("H36").Style = "Normal"
$ws.Range("H36").HorizontalAlignment = -4108
$ws.Range("H36").VerticalAlignment = -4108
$ws.Range("H36").WrapText = $true
$ws.Range("H36").Value = "Y"

$ws.Range("H37").Style = "Normal"
$ws.Range("H37").HorizontalAlignment = -4108
$ws.Range("H37").VerticalAlignment = -4108
$ws.Range("H37").WrapText = $true
$ws.Range("H37").Value = "Y"

$ws.Range("H38").Style = "Normal"
$ws.Range("H38").HorizontalAlignment = -4108
$ws.Range("H38").VerticalAlignment = -4108
$ws.Range("H38").WrapText = $true
$ws.Range("H38").Value = "Y"

# ---------------------------------------------------------------------------
# 4. Hyperlinks: row-insert shifted the cell contents but not the hyperlink
#    anchors, so re-anchor the three manifest hyperlinks that used to live
#    on E25/E26/E27 onto their new homes E36/E37/E38.
# ---------------------------------------------------------------------------
$moves = @{ 'E25' = 'E36'; 'E26:E27' = 'E36:E38'; 'E26' = 'E37'; 'E27' = 'E38' }

$old1 = $null
$old2 = $null
$old3 = $null
$old4 = $null
foreach ($h in @($ws.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$E$25') { $old1 = @{ Address = $h.Address; Tip = $h.ScreenTip } }
    elseif ($addr -eq '$E$26:$E$27') { $old2 = @{ Address = $h.Address; Display = $h.TextToDisplay; Tip = $h.ScreenTip } }
    elseif ($addr -eq '$E$26') { $old3 = @{ Address = $h.Address; Tip = $h.ScreenTip } }
    elseif ($addr -eq '$E$27') { $old4 = @{ Address = $h.Address; Tip = $h.ScreenTip } }
}
foreach ($h in @($ws.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$E$25' -or $addr -eq '$E$26:$E$27' -or $addr -eq '$E$26' -or $addr -eq '$E$27') {
        $h.Delete()
    }
}
if ($old1) { $ws.Hyperlinks.Add($ws.Range("E36"), $old1.Address) | Out-Null }
if ($old2) { $ws.Hyperlinks.Add($ws.Range("E36:E38"), $old2.Address, "", "", $old2.Display) | Out-Null }
if ($old3) { $ws.Hyperlinks.Add($ws.Range("E37"), $old3.Address) | Out-Null }
if ($old4) { $ws.Hyperlinks.Add($ws.Range("E38"), $old4.Address) | Out-Null }

# ---------------------------------------------------------------------------
# 5. View state: the active selection moved to H39, and the sheet is
#    scrolled so row 12 is at the top of the frozen pane.
# ---------------------------------------------------------------------------
$ws.Range("H39").Select()
$excel.ActiveWindow.ScrollRow = 12
